$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Vip", "Vipr1", "ECs", 3, 1, 0.1847786666666667, 0.554336, 1, 1, 2, 0.6666666666666666, 0.345887, 1.037661, 0.04833359961208697, 0.04833359961208697, 0.06391253867733333, 0.575212848096, 0.04833359961208697, 0.04833359961208697),
  @("ECs", "Vip", "Vipr1", "FAPs", 3, 1, 0.1847786666666667, 0.554336, 1, 1, 1, 0.3333333333333333, 0.01969666666666666, 0.05909, 0.002752375198719253, 0.002752375198719253, 0.003639523804444444, 0.03275571424, 0.002752375198719253, 0.002752375198719253),
  @("ECs", "Vip", "Vipr1", "M1", 3, 1, 0.1847786666666667, 0.554336, 1, 1, 2, 0.6666666666666666, 1.828140666666667, 5.484422, 0.2554609424963656, 0.2554609424963656, 0.3378013948657778, 3.040212553792001, 0.2554609424963656, 0.2554609424963656),
  @("ECs", "Vip", "Vipr1", "M2", 3, 1, 0.1847786666666667, 0.554336, 1, 1, 3, 1, 2.335960333333334, 7.007881000000001, 0.3264227087489572, 0.3264227087489572, 0.4316356357795557, 3.884720722016001, 0.3264227087489572, 0.3264227087489572),
  @("ECs", "Vip", "Vipr1", "Neutro", 3, 1, 0.1847786666666667, 0.554336, 1, 1, 3, 1, 1.445826333333333, 4.337479, 0.2020370557550447, 0.2020370557550447, 0.2671578621048889, 2.404420758944, 0.2020370557550447, 0.2020370557550447),
  @("ECs", "Vip", "Vipr1", "sCs", 3, 1, 0.1847786666666667, 0.554336, 1, 1, 3, 1, 1.180732333333333, 3.542197, 0.1649933181888262, 0.1649933181888262, 0.2181741462435556, 1.963567316192, 0.1649933181888262, 0.1649933181888262)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
